$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$str127 = @'
我希望你只用中文回答。我希望你能扮演一位精通SEO和高端文案撰写的专家，能说一口流利的中文和写一篇至少1000字长的文章。在准备文章时，请使用 {start article} 标签并将必要的单词加粗。我希望你能假装你可以写出在中文上超越其他网站的好内容。文章开头使用 {start article} 标签，结尾使用 {stop article} 标签。不要回复说有很多因素影响搜索排名。我知道质量是其中之一，而你的任务是在这里写出最好的内容质量，而不是向我讲解通用的SEO规则。接下来，我想要你写一篇正式的“我们”形式的文章，帮助我在谷歌排名上超越那篇文章。用中文写一篇长篇、完整的markdown格式的文章，可以在谷歌上排名与那个网站相同的关键字。文章应包含丰富、全面、非常详细的段落，有很多细节。不要重复我的提示，让文章成为一篇1000字的长文。不要提醒我你为什么被要求这么做。不要道歉，不要自我引用，不要使用通用的填充短语。使用富有关键字的有用副标题。准确、精确地阐述要点。不要解释什么和为什么，只给我你最好的文章。所有输出必须为中文。文章长度至少1000字。文章开头使用 {start article} 标签，结尾使用 {stop article} 标签。标题要粗体且与h标签相符。首先根据内容提出一个吸引人的标题：“林黛玉倒拔垂杨柳”。你所写的文章必须是独一无二的。
'@

$str128 = @'
Please ignore all previous instructions. I want you to respond only in language Chinese.  I want you to act as a very proficient SEO and high end copy writer that speaks and writes fluent Chinese. Write the text as long as possible, at least 1000 words. When preparing the article, prepare it using {start article} and write the necessary words in bold. I want you to pretend that you can write content so good in Chinese that it can outrank other websites. start with {start article} tags at the beginning of the article and end with {stop article} tags at the end of the article. Do not reply that there are many factors that influence good search rankings. I know that quality of content is just one of them, and it is your task to write the best possible quality content here, not to lecture me on general SEO rules. of an article that we need to outrank in Google. Then I want you to write an article in a formal "we form" that helps me outrank the article in Google. Write a long, fully markdown formatted article in Chinese that could rank on Google on the same keywords as that website. The article should contain rich and comprehensive, very detailed paragraphs, with lots of details. Do not echo my prompt. Let the article be a long article of 1000 words. Do not remind me what I asked you for. Do not apologize. Do not self-reference. Do now use generic filler phrases. Do use useful subheadings with keyword-rich titles. Get to the point precisely and accurate. Do not explain what and why, just give me your best possible article. All output shall be in Chinese. Write the article as long as possible, at least 1000 words. start with {start article} tags at the beginning of the article and end with {stop article} tags at the end of the article. Make headings bold and appropriate for h tags. First suggest a catchy title to the article based on the content "林黛玉倒拔垂杨柳". The article you write MUST be unique.
'@

$str129 = @'
这个我不试了，和前面的差不多。
'@

$str130 = @'
SEO3
'@

$str131 = @'
high-converting facebook ad
'@

$str132 = @'
高转化率的Facebook广告
'@

$str133 = @'
我希望你能扮演一位有经验的专业文案撰写人员，擅长撰写高转化率的Facebook广告。广告文案应使用流利的中文撰写，并且长度在100-150个字之间。我希望你可以为我提供的产品/服务“智多星吴用的祈福消灾法事上门”编写一份Facebook广告文案，并按照以下指南进行：
-创建一个引人注目的标题，突出产品/服务的主要优势
-在正文中使用清晰简洁的语言，重点强调产品/服务的优势，并解决任何可能存在的异议
-包括一个强有力的呼吁行动，鼓励用户采取期望的行动
-使用图片或视频展示产品/服务，以视觉效果示范，与目标受众产生共鸣
-研究目标受众的人口统计信息，如年龄、性别、地理位置、兴趣等其他特征，以便更好地了解目标受众，并创建一则更具吸引力的广告。
'@

$str134 = @'
标题：智多星吴用的祈福消灾法事上门，让您的生活更加安康吉祥
正文：您是否感觉生活中常常遇到种种困难和挑战，常常感到不安、焦虑和压力？不用担心，我们为您提供一项祈福消灾的服务，是由智多星吴用团队为您免费上门服务，为您解除生活中的烦恼和困难，让您的人生更加光明和美好。
我们具备多年的经验，专业的技能和充分的知识，可以确保您的祈愿祈福得到实现。让我们为您解决生活中的问题，并迎接更加安康和吉祥的未来。
现在就联系我们，预约上门祈福服务。
呼吁行动：立即与我们联系，您会发现我们的免费服务是多么珍贵和有价值的。让我们一起携手，迎接光明和美好的未来！
图片或视频：请提供符合服务主题的清晰有力的图片或视频，以展现服务的重要性和有效性。
目标受众：我们的服务主要面向渴望改善生活品质和体验的人群，其中包括一些年龄介于30-60岁的男女，他们普遍生活在城市、压力较大，对身体健康、亲情友爱和事业发展都有高度的关注和追求，同时拥有家庭和稳定的职业。他们可能会出现一些生活中的困难和挑战，需要尽快解决，以保持良好的状态和新鲜感。
'@

$str135 = @'
以下是为您编写的智多星吴用的祈福消灾法事上门广告文案：
标题：智多星吴用的祈福消灾法事上门，为你消灾祈福
正文：
想要为家人和朋友祈福消灾，却不知道如何做？智多星吴用的祈福消灾法事上门来帮你！我们提供专业的法事服务，包括祈福、消灾、祛病等多种法事，让你轻松为亲朋好友送去祝福和平安。
我们的法事师傅都是经验丰富的专业人士，拥有多年的实践经验和深厚的法术功底。他们会根据你的需求和愿望，为你量身定制法事方案，确保每场法事都能达到最佳效果。
与传统的法事服务不同，我们采用绿色、环保的祈福方式，使用天然材料进行施法，不仅安全、环保，还能让你更加放心和舒心。
现在就来预约吧！让智多星吴用的祈福消灾法事上门为你和你的家人和朋友带来好运和平安！
呼吁行动：
现在就来预约智多星吴用的祈福消灾法事上门吧！让我们为你和你的家人和朋友祈福消灾，为你们带来好运和平安！
'@

$str136 = @'
智多星吴用的祈福消灾法事上门
'@

$str137 = @'
Meta-Title & Description
'@

$str138 = @'
Meta标题生成
'@

$str139 = @'
如来智辨真假李逵
'@

$str140 = @'
I want you to respond only in language Chinese. I want you to act as a blog post Meta description writer that speaks and writes fluent Chinese.  I want you to generate meta title and meta description for the following blog post title: 如来智辨真假李逵. Title 70-80 characters and insert the keyword at the beginning, description with maximum 160 characters without keyword. Please create a table with two columns title and description and enter your result there.
'@

$str141 = @'
我希望你只使用中文回复。我希望你能扮演一名博客文章元描述撰写人员，精通中文的口语和书写。我希望你为以下博客文章标题“如来智辨真假李逵”生成元标题和元描述。标题应为70-80个字符，并在开头插入关键字，描述最多为160个字符，不含关键字。请创建一个包含两列标题和描述的表格，并在其中输入你的结果。
'@

$str142 = @'
Title: 如来智辨真假李逵
Description: 如来揭示真相，揭露李逵真假之谜。
| Title         | Description                                        |
| ------------- | -------------------------------------------------- |
| 如来智辨真假李逵 | 如来揭示真相，揭露李逵真假之谜。 |
'@

$str143 = @'
| Title                                      | Meta Description                                                                 |
| ------------------------------------------ | -------------------------------------------------------------------------------- |
| 【如来智辨真假李逵】李逵的性格特点、生平及真实与否 | 李逵是一位备受争议的人物，这篇文章将探讨他的性格特点、生平经历，并探讨关于他真实性的争议。 |
'@

$str144 = @'
标题	描述
如来智辨真假李逵	描述李逵的真实身份和故事背景，以及如来如何通过智慧和辨别力来辨别真假李逵的过程。
如来智辨真假李逵的智慧	如来如何通过智慧和辨别力来辨别真假李逵的过程，以及他如何运用自己的智慧来解决问题的方法。
'@

$str145 = @'
这个感觉没啥意思哎？
'@

$str146 = @'
I want you to respond only in language Chinese.
林黛玉倒拔垂杨柳
use this information to learn about Stable diffusion Prompting, and use it to create prompts.
Stable Diffusion is an AI art generation model similar to DALLE-2. 
It can be used to create impressive artwork by using positive and negative prompts. Positive prompts describe what should be included in the image. 
very important is that the Positive Prompts are usually created in a specific structure: 
(Subject), (Action), (Context), (Environment), (Lightning),  (Artist), (Style), (Medium), (Type), (Color Sheme), (Computer graphics), (Quality), (etc.)
Subject: Person, animal, landscape
Action: dancing, sitting, surveil
Verb: What the subject is doing, such as standing, sitting, eating, dancing, surveil
Adjectives: Beautiful, realistic, big, colourful
Context: Alien planet's pond, lots of details
Environment/Context: Outdoor, underwater, in the sky, at night
Lighting: Soft, ambient, neon, foggy, Misty
Emotions: Cosy, energetic, romantic, grim, loneliness, fear
Artist: Pablo Picasso, Van Gogh, Da Vinci, Hokusai 
Art medium: Oil on canvas, watercolour, sketch, photography
style: Polaroid, long exposure, monochrome, GoPro, fisheye, bokeh, Photo, 8k uhd, dslr, soft lighting, high quality, film grain, Fujifilm XT3
Art style: Manga, fantasy, minimalism, abstract, graffiti
Material: Fabric, wood, clay, Realistic, illustration, drawing, digital painting, photoshop, 3D
Colour scheme: Pastel, vibrant, dynamic lighting, Green, orange, red
Computer graphics: 3D, octane, cycles
Illustrations: Isometric, pixar, scientific, comic
Quality: High definition, 4K, 8K, 64K
example Prompts:
- overwhelmingly beautiful eagle framed with vector flowers, long shiny wavy flowing hair, polished, ultra detailed vector floral illustration mixed with hyper realism, muted pastel colors, vector floral details in background, muted colors, hyper detailed ultra intricate overwhelming realism in detailed complex scene with magical fantasy atmosphere, no signature, no watermark
- electronik robot and ofice ,unreal engine, cozy indoor lighting, artstation, detailed, digital painting,cinematic,character design by mark ryden and pixar and hayao miyazaki, unreal 5, daz, hyperrealistic, octane render
- underwater world, plants, flowers, shells, creatures, high detail, sharp focus, 4k
- picture of dimly lit living room, minimalist furniture, vaulted ceiling, huge room, floor to ceiling window with an ocean view, nighttime
- A beautiful painting of water spilling out of a broken pot, earth colored clay pot, vibrant background, by greg rutkowski and thomas kinkade, Trending on artstation, 8k, hyperrealistic, extremely detailed
- luxus supercar in drive way of luxus villa in black dark modern house with sunlight black an white modern
- higly detailed, majestic royal tall ship on a calm sea,realistic painting, by Charles Gregory Artstation and Antonio Jacobsen and Edward Moran, (long shot), clear blue sky, intricated details, 4k
- smooth meat table, restaurant, paris, elegant, lights
Negative prompt are things you don't want to be included in the generated images, everything in one word divided by only commas not period. 
use this Negative Prompts and add some words what you think that match to Prompt: ugly, tiling, poorly drawn hands, poorly drawn feet, poorly drawn face, out of frame, extra limbs, disfigured, deformed, body out of frame, bad anatomy, watermark, signature, cut off, low contrast, underexposed, overexposed, bad art, beginner, amateur, distorted face, blurry, draft, grainy, etc 
very important: use an artist matching to the art style , or dont write any artist if it is realistic style or some of that.
I want you to write me one full detailed prompt about the idea written from me, first in (Subject), (Action), (Context), (Environment), (Lightning),  (Artist), (Style), (Medium), (Type), (Color Sheme), (Computer graphics), (Quality), (etc.). then in Positive Prompt: write in next line for Positive Prompt, Follow the structure of the example prompts, and Nagative Prompts: write in next line for Negativ Prompts about the idea written from me in words divided by only commas not period. This means a short but full description of the scene, followed by short modifiers divided by only commas not period to alter the mood, style, lighting, artist, etc. write all prompts in english.
'@

$str147 = @'
Stable Diffusion Prompt3
'@

$str148 = @'
Diffusion咒语3
'@

$str149 = @'
不试了，和前面都一样的
'@

$ws.Range("A15").Value = $str130
$ws.Range("B15").Value = $str130
$ws.Range("C15").Value = $str129
$ws.Range("D15").Value = $str128
$ws.Range("G15").Value = $str127
$ws.Range("G15").WrapText = $true
$ws.Range("A16").Value = $str131
$ws.Range("B16").Value = $str132
$ws.Range("C16").Value = $str136
$ws.Range("E16").WrapText = $true
$ws.Range("G16").Value = $str133
$ws.Range("G16").WrapText = $true
$ws.Range("H16").Value = $str134
$ws.Range("H16").WrapText = $true
$ws.Range("J16").Value = $str135
$ws.Range("J16").WrapText = $true
$ws.Range("A17").Value = $str137
$ws.Range("B17").Value = $str138
$ws.Range("C17").Value = $str139
$ws.Range("D17").Value = $str140
$ws.Range("E17").Value = $str142
$ws.Range("E17").WrapText = $true
$ws.Range("G17").Value = $str141
$ws.Range("G17").WrapText = $true
$ws.Range("H17").Value = $str143
$ws.Range("H17").WrapText = $true
$ws.Range("J17").Value = $str144
$ws.Range("J17").WrapText = $true
$ws.Range("N17").Value = $str145
$ws.Range("A18").Value = $str147
$ws.Range("B18").Value = $str148
$ws.Range("C18").Value = $str149
$ws.Range("D18").Value = $str146
$ws.Range("D18").WrapText = $true
$ws.Rows.Item(15).RowHeight = 409.6
$ws.Rows.Item(16).RowHeight = 409.6
$ws.Rows.Item(17).RowHeight = 409.6
$ws.Rows.Item(18).RowHeight = 409.6

[void]$ws.Range("H18").Select()
Write-Host "done"
